# Updated cryptos list (price/volume refresh) on Sat Sep  7 06:00:22 UTC 2024
# Forces numeric-looking Price values to remain text (matching the source
# data's inline-string typing) by using a leading apostrophe, then clears
# the resulting quote-prefix style so no stray per-cell formatting is left
# behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.026.43"
$ws.Range("E2").Value = "  -4.14%  "
$ws.Range("D3").Value = "2.261.28"
$ws.Range("E3").Value = "  -4.74%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'491.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.09%  "
$ws.Range("D6").Value = "'126.88"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.61%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("E8").Value = "  -3.37%  "
$ws.Range("D9").Value = "2.260.70"
$ws.Range("E9").Value = "  -5.13%  "
$ws.Range("D10").Value = "'0.0931"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.80%  "
$ws.Range("D11").Value = "'0.150"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.09%  "
$ws.Range("D12").Value = "'0.323"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.61%  "
$ws.Range("D13").Value = "'4.63"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.31%  "
$ws.Range("D14").Value = "2.660.39"
$ws.Range("E14").Value = "  -4.82%  "
$ws.Range("D15").Value = "'21.44"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.53%  "
$ws.Range("D16").Value = "54.018.47"
$ws.Range("E16").Value = "  -4.11%  "
$ws.Range("E17").Value = "  -3.29%  "
$ws.Range("D18").Value = "2.278.80"
$ws.Range("E18").Value = "  -4.20%  "
$ws.Range("D19").Value = "'9.77"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.75%  "
$ws.Range("E20").Value = "  -0.33%  "
$ws.Range("D21").Value = "'297.61"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.99%  "
$ws.Range("E22").Value = "  -0.89%  "
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("D24").Value = "'63.81"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.15%  "
$ws.Range("D25").Value = "'1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.22%  "
$ws.Range("E26").Value = "  +0.31%  "
$ws.Range("D27").Value = "2.363.95"
$ws.Range("E27").Value = "  -4.78%  "
$ws.Range("E28").Value = "  -1.18%  "
$ws.Range("D29").Value = "'7.11"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.53%  "
$ws.Range("D30").Value = "'162.74"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.89%  "
$ws.Range("E31").Value = "  -3.23%  "
$ws.Range("E32").Value = "  -4.96%  "
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("D34").Value = "'5.79"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.42%  "
$ws.Range("D35").Value = "'0.997"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.11%  "
$ws.Range("E36").Value = "  -0.90%  "
$ws.Range("D37").Value = "'17.41"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.25%  "
$ws.Range("E38").Value = "  +0.31%  "
$ws.Range("D39").Value = "'0.838"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.62%  "
$ws.Range("D40").Value = "'3.62"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.07%  "
$ws.Range("D41").Value = "'35.37"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.82%  "
$ws.Range("E42").Value = "  +0.84%  "
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("E44").Value = "  -1.64%  "
$ws.Range("D45").Value = "'126.58"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.61%  "
$ws.Range("D46").Value = "'4.79"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.32%  "
$ws.Range("E47").Value = "  -0.89%  "
$ws.Range("E48").Value = "  +1.77%  "
$ws.Range("E49").Value = "  -3.76%  "
$ws.Range("E50").Value = "  -0.91%  "
$ws.Range("E51").Value = "  -2.13%  "
